$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy formatting from row 37 down to the two new rows (38 and 39)
# so the new rows match the style of the existing data rows.
$ws.Range("A37:E37").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)
$ws.Range("A37:E37").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update existing row 37, Results column (E) from PASS to SKIP
$ws.Range("E37").Value = "SKIP"

# Add new row 38: CommentsTabTimeStampValidationTest
$ws.Range("A38").Value = "CommentsTabTimeStampValidationTest"
$ws.Range("B38").Value = "TBD"
$ws.Range("C38").Value = "Verify that Comments tab comments displayed with timestamp"
$ws.Range("D38").Value = "Y"
$ws.Range("E38").Value = "SKIP"

# Add new row 39: HCRProfileBadgeTest
$ws.Range("A39").Value = "HCRProfileBadgeTest"
$ws.Range("B39").Value = "TBD"
$ws.Range("C39").Value = "Verify that HCR profile having badge along with their name"
$ws.Range("D39").Value = "Y"
$ws.Range("E39").Value = "PASS"

# Update sheet view: scroll so row 13 is at the top, and select C34
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C34").Select()
